# The commit deletes the original first slide ("Data Exploration Challenge" /
# "Group Members" title slide). Every other slide's content is unchanged -
# they simply shift up one position (what was slide 2 becomes slide 1, etc.)
# as a natural consequence of removing slide 1. The presentation's sldIdLst
# simply loses its first entry (id=256) while every other <p:sldId> keeps
# its original id and relative order.

$p = $ppt.ActivePresentation
$p.Slides.Item(1).Delete()
